$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab
$ws.Name = "SCD0018"

# 2) Update the TC_ID column (B2:B8) from "DGS-317" to "SCD0018-025"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 2).Value = "SCD0018-025"
}

# 3) Re-apply horizontal=left / vertical=center alignment across the whole used range
$rng = $ws.Range("A1:V8")
$rng.HorizontalAlignment = -4131
$rng.VerticalAlignment = -4108

# 4) Column B needs to widen to fit the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 12.28515625

# 5) Restore the view/selection state (scrolled so row 8 is on top, B9 selected)
$ws.Range("B9").Select()
$excel.ActiveWindow.ScrollRow = 8
